# Update the "Sheet1" data for fiscal year 2025 (row 7) with refreshed figures
# ("atualizei dados bibi e add").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2739981.91
$ws.Range("C7").Value = -38.33144869395184
$ws.Range("D7").Value = 2789
$ws.Range("E7").Value = 2789
$ws.Range("F7").Value = 982.4244926496953
$ws.Range("G7").Value = 4.71934707258661
